$wb = $excel.ActiveWorkbook

# Update the status text on all three sheets where "Ready for handoff" appears
# ("Overview" columns E/F, and the "Status" column (C) on the per-locale
# "zh-cn"/"de-de" sheets), then shrink the now-narrower status columns to
# match the shorter "In Translation" text (mirrors the width change Excel's
# own AutoFit produces when the cell text gets shorter).
#
# Target stored column width from the real-Excel AutoFit is 13.4101845877511
# (down from 17.2159881591797). This runtime quantizes ColumnWidth to a 6px
# per-character grid, so 12.5 is the input that lands closest to that target.
$newWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = $newWidth
